# Weekly CompStat data refresh: report number, date range, and updated crime counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: report volume/number and week-covering date range ---
$ws.Range("A8").Value = "Volume 32   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Crime Complaints + Transit/Housing tables: weekly/28-day/YTD counts and % changes ---
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 37
$ws.Range("K16").Value = 48
$ws.Range("L16").Value = -15.909090909090
$ws.Range("M16").Value = -32.727272727272
$ws.Range("N16").Value = -81.122448979591
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 62
$ws.Range("J17").Value = 62
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -16.216216216216
$ws.Range("M17").Value = 31.914893617021
$ws.Range("N17").Value = -17.333333333333
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = -15.384615384615
$ws.Range("L18").Value = -38.888888888888
$ws.Range("M18").Value = -85.897435897435
$ws.Range("N18").Value = -94.811320754717
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 125
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 7.692307692307
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 105
$ws.Range("K19").Value = -15.238095238095
$ws.Range("L19").Value = 23.611111111111
$ws.Range("M19").Value = 34.848484848484
$ws.Range("N19").Value = -23.275862068965
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 15
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = -34.782608695652
$ws.Range("L20").Value = -43.396226415094
$ws.Range("M20").Value = -16.666666666666
$ws.Range("N20").Value = -95.810055865921
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 23.529411764705
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -9.230769230769
$ws.Range("I21").Value = 235
$ws.Range("J21").Value = 256
$ws.Range("K21").Value = -8.203125
$ws.Range("L21").Value = -11.654135338345
$ws.Range("M21").Value = -17.253521126760
$ws.Range("N21").Value = -82.264150943396
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = -100
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = -83.333333333333
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -27.272727272727
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = -2.857142857142
$ws.Range("L23").Value = -5.555555555555
$ws.Range("M23").Value = 142.857142857143
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -57.142857142857
$ws.Range("F24").Value = 31
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = -24.390243902439
$ws.Range("I24").Value = 166
$ws.Range("J24").Value = 166
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -19.417475728155
$ws.Range("M24").Value = 13.698630136986
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 60
$ws.Range("I25").Value = 41
$ws.Range("J25").Value = 34
$ws.Range("K25").Value = 20.588235294117
$ws.Range("L25").Value = -21.153846153846
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -54.545454545454
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -7.692307692307
$ws.Range("I26").Value = 106
$ws.Range("J26").Value = 104
$ws.Range("K26").Value = 1.923076923076
$ws.Range("L26").Value = -7.826086956521
$ws.Range("M26").Value = -15.873015873015
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 166.666666666667
$ws.Range("L27").Value = 166.666666666667
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 57.142857142857
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C29").Copy($ws.Range("F29"))
$ws.Range("G29").Value = 1
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = -100
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = -40
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C30").Copy($ws.Range("F30"))
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -40
